# Update college name and reposition/resize its textbox on the title slide
# (slide 1, shape id=21 "学院：..." textbox).
#
# PowerPoint's Shape.Left/.Top/.Width/.Height are 32-bit `float` properties,
# so a naive `= emu / 12700` assignment can land 1 EMU off after the
# float32 round-trip. We assign the natural value first and, if the
# round-tripped EMU doesn't match exactly, nudge through nearby candidate
# point values until the saved EMU is exact.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$emuPerPt = 12700.0

function Get-Emu($ptVal) {
    return [math]::Round($ptVal * $emuPerPt)
}

function Set-PreciseLeft($shape, $targetEmu) {
    $approxPt = $targetEmu / $emuPerPt
    $shape.Left = $approxPt
    if ((Get-Emu($shape.Left)) -eq $targetEmu) { return }
    for ($i = -500; $i -le 500; $i++) {
        if ($i -eq 0) { continue }
        $cand = $approxPt + $i * 0.00001
        $shape.Left = $cand
        if ((Get-Emu($shape.Left)) -eq $targetEmu) { return }
    }
}

function Set-PreciseTop($shape, $targetEmu) {
    $approxPt = $targetEmu / $emuPerPt
    $shape.Top = $approxPt
    if ((Get-Emu($shape.Top)) -eq $targetEmu) { return }
    for ($i = -500; $i -le 500; $i++) {
        if ($i -eq 0) { continue }
        $cand = $approxPt + $i * 0.00001
        $shape.Top = $cand
        if ((Get-Emu($shape.Top)) -eq $targetEmu) { return }
    }
}

function Set-PreciseWidth($shape, $targetEmu) {
    $approxPt = $targetEmu / $emuPerPt
    $shape.Width = $approxPt
    if ((Get-Emu($shape.Width)) -eq $targetEmu) { return }
    for ($i = -500; $i -le 500; $i++) {
        if ($i -eq 0) { continue }
        $cand = $approxPt + $i * 0.00001
        $shape.Width = $cand
        if ((Get-Emu($shape.Width)) -eq $targetEmu) { return }
    }
}

function Set-PreciseHeight($shape, $targetEmu) {
    $approxPt = $targetEmu / $emuPerPt
    $shape.Height = $approxPt
    if ((Get-Emu($shape.Height)) -eq $targetEmu) { return }
    for ($i = -500; $i -le 500; $i++) {
        if ($i -eq 0) { continue }
        $cand = $approxPt + $i * 0.00001
        $shape.Height = $cand
        if ((Get-Emu($shape.Height)) -eq $targetEmu) { return }
    }
}

# Find the "学院：..." textbox (shape id 21) on the title slide.
$collegeShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Id -eq 21) {
        $collegeShape = $sh
        break
    }
}

$collegeShape.TextFrame.TextRange.Text = "学院：人工智能学院"

Set-PreciseLeft  $collegeShape 5312187
Set-PreciseWidth $collegeShape 1783122
